$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (20) loses the "latest row" date format and reverts
# to the standard date/time format used by all other data rows.
$ws.Range("A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 21.
$ws.Range("A21").Value = 45761
$ws.Range("A21").NumberFormat = "YYYY-MM-DD"

$ws.Range("B21").Value = 82
$ws.Range("C21").Value = 82
$ws.Range("D21").Value = 80
